$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 42.75280866666666
$ws.Range("H2").Value = 128.258426
$ws.Range("I2").Value = 0.8529286054750734
$ws.Range("J2").Value = 0.8529286054750735
$ws.Range("M2").Value = 0.06743766666666666
$ws.Range("N2").Value = 0.202313
$ws.Range("O2").Value = 0.004349811883262162
$ws.Range("P2").Value = 0.004349811883262163
$ws.Range("Q2").Value = 2.883149659926444
$ws.Range("R2").Value = 25.94834693933799
$ws.Range("S2").Value = 0.003710078983669699
$ws.Range("T2").Value = 0.0037100789836697
$ws.Range("G3").Value = 42.75280866666666
$ws.Range("H3").Value = 128.258426
$ws.Range("I3").Value = 0.8529286054750734
$ws.Range("J3").Value = 0.8529286054750735
$ws.Range("O3").Value = 0.7938207485680675
$ws.Range("P3").Value = 0.7938207485680676
$ws.Range("Q3").Value = 526.1616094441662
$ws.Range("R3").Value = 4735.454484997495
$ws.Range("S3").Value = 0.6770724240733407
$ws.Range("T3").Value = 0.6770724240733409
$ws.Range("G4").Value = 42.75280866666666
$ws.Range("H4").Value = 128.258426
$ws.Range("I4").Value = 0.8529286054750734
$ws.Range("J4").Value = 0.8529286054750735
$ws.Range("O4").Value = 0.2018294395486703
$ws.Range("P4").Value = 0.2018294395486704
$ws.Range("Q4").Value = 133.7769300408209
$ws.Range("R4").Value = 1203.992370367388
$ws.Range("S4").Value = 0.172146102418063
$ws.Range("T4").Value = 0.1721461024180631
$ws.Range("I5").Value = 0.04642608686423023
$ws.Range("J5").Value = 0.04642608686423023
$ws.Range("M5").Value = 0.06743766666666666
$ws.Range("N5").Value = 0.202313
$ws.Range("O5").Value = 0.004349811883262162
$ws.Range("P5").Value = 0.004349811883262163
$ws.Range("Q5").Value = 0.1569338344324444
$ws.Range("R5").Value = 1.412404509892
$ws.Range("S5").Value = 0.00020194474433539
$ws.Range("T5").Value = 0.0002019447443353901
$ws.Range("I6").Value = 0.04642608686423023
$ws.Range("J6").Value = 0.04642608686423023
$ws.Range("O6").Value = 0.7938207485680675
$ws.Range("P6").Value = 0.7938207485680676
$ws.Range("S6").Value = 0.03685399102764936
$ws.Range("T6").Value = 0.03685399102764937
$ws.Range("I7").Value = 0.04642608686423023
$ws.Range("J7").Value = 0.04642608686423023
$ws.Range("O7").Value = 0.2018294395486703
$ws.Range("P7").Value = 0.2018294395486704
$ws.Range("S7").Value = 0.009370151092245472
$ws.Range("T7").Value = 0.009370151092245474
$ws.Range("G8").Value = 5.044817999999999
$ws.Range("I8").Value = 0.1006453076606963
$ws.Range("J8").Value = 0.1006453076606963
$ws.Range("M8").Value = 0.06743766666666666
$ws.Range("N8").Value = 0.202313
$ws.Range("O8").Value = 0.004349811883262162
$ws.Range("P8").Value = 0.004349811883262163
$ws.Range("Q8").Value = 0.3402107546779999
$ws.Range("R8").Value = 3.061896792101999
$ws.Range("S8").Value = 0.0004377881552570731
$ws.Range("T8").Value = 0.0004377881552570732
$ws.Range("G9").Value = 5.044817999999999
$ws.Range("I9").Value = 0.1006453076606963
$ws.Range("J9").Value = 0.1006453076606963
$ws.Range("O9").Value = 0.7938207485680675
$ws.Range("P9").Value = 0.7938207485680676
$ws.Range("R9").Value = 558.782142486984
$ws.Range("S9").Value = 0.07989433346707739
$ws.Range("T9").Value = 0.0798943334670774
$ws.Range("G10").Value = 5.044817999999999
$ws.Range("I10").Value = 0.1006453076606963
$ws.Range("J10").Value = 0.1006453076606963
$ws.Range("O10").Value = 0.2018294395486703
$ws.Range("P10").Value = 0.2018294395486704
$ws.Range("S10").Value = 0.02031318603836183
$ws.Range("T10").Value = 0.02031318603836183
